# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-17 (serial 45186) to 2023-09-19 (serial 45188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
